$wb = $excel.ActiveWorkbook

# === Sheet: Tarantula ===
$ws = $wb.Worksheets.Item("Tarantula")
$ws.Cells.Item(2, 3).Value = -1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 3).Value = -1
$ws.Cells.Item(3, 4).Value = 25
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = 16
$ws.Cells.Item(4, 3).Value = -1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 3).Value = -1
$ws.Cells.Item(5, 4).Value = 8
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 9).Value = 6
$ws.Cells.Item(6, 3).Value = -1
$ws.Cells.Item(6, 4).Value = 11
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 9).Value = 10
$ws.Cells.Item(7, 3).Value = -1
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 9).Value = 28
$ws.Cells.Item(8, 3).Value = -1
$ws.Cells.Item(8, 4).Value = 27
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 9).Value = 13
$ws.Cells.Item(9, 3).Value = -1
$ws.Cells.Item(9, 4).Value = 27
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 9).Value = 13
$ws.Cells.Item(10, 3).Value = -1
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 9).Value = 11
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(12, 3).Value = -1
$ws.Cells.Item(12, 4).Value = 9
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 3).Value = -1
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 9).Value = 6
$ws.Cells.Item(14, 3).Value = -1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 3).Value = -1
$ws.Cells.Item(15, 4).Value = 8
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 9).Value = 6
$ws.Cells.Item(18, 3).Value = -1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(19, 3).Value = -1
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 3).Value = -1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 3).Value = -1
$ws.Cells.Item(21, 4).Value = 11
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 9).Value = 11
$ws.Cells.Item(22, 3).Value = -1
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 3).Value = -1
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(24, 3).Value = -1
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(25, 3).Value = -1
$ws.Cells.Item(25, 4).Value = 26
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 9).Value = 13
$ws.Cells.Item(26, 3).Value = -1
$ws.Cells.Item(26, 4).Value = 6
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = -1
$ws.Cells.Item(27, 4).Value = 26
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 9).Value = 13
$ws.Cells.Item(28, 3).Value = -1
$ws.Cells.Item(28, 4).Value = 22
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = 10
$ws.Cells.Item(29, 3).Value = -1
$ws.Cells.Item(29, 4).Value = 8
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = 12
$ws.Cells.Item(30, 3).Value = -1
$ws.Cells.Item(30, 4).Value = 5
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 3).Value = -1
$ws.Cells.Item(31, 4).Value = 8
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 3).Value = -1
$ws.Cells.Item(32, 4).Value = 8
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 9).Value = 7
$ws.Cells.Item(33, 3).Value = -1
$ws.Cells.Item(33, 4).Value = 9
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(34, 3).Value = -1
$ws.Cells.Item(34, 4).Value = 40
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 9).Value = 13

# === Sheet: Ochiai ===
$ws = $wb.Worksheets.Item("Ochiai")
$ws.Cells.Item(2, 3).Value = -1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 3).Value = -1
$ws.Cells.Item(3, 4).Value = 5
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(4, 3).Value = -1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 3).Value = -1
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 3).Value = -1
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(7, 3).Value = -1
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 3).Value = -1
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 9).Value = 3
$ws.Cells.Item(9, 3).Value = -1
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 9).Value = 3
$ws.Cells.Item(10, 3).Value = -1
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 9).Value = 2
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(12, 3).Value = -1
$ws.Cells.Item(12, 4).Value = 9
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 3).Value = -1
$ws.Cells.Item(13, 4).Value = 6
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(14, 3).Value = -1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 3).Value = -1
$ws.Cells.Item(15, 4).Value = 6
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(18, 3).Value = -1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(19, 3).Value = -1
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 3).Value = -1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 3).Value = -1
$ws.Cells.Item(21, 4).Value = 5
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 9).Value = 2
$ws.Cells.Item(22, 3).Value = -1
$ws.Cells.Item(22, 4).Value = 14
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 3).Value = -1
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(24, 3).Value = -1
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(25, 3).Value = -1
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(26, 3).Value = -1
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = -1
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 9).Value = 3
$ws.Cells.Item(28, 3).Value = -1
$ws.Cells.Item(28, 4).Value = 18
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = 9
$ws.Cells.Item(29, 3).Value = -1
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = 2
$ws.Cells.Item(30, 3).Value = -1
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 3).Value = -1
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 3).Value = -1
$ws.Cells.Item(32, 4).Value = 6
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 9).Value = 4
$ws.Cells.Item(33, 3).Value = -1
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(34, 3).Value = -1
$ws.Cells.Item(34, 4).Value = 9
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 9).Value = 3

# === Sheet: Op2 ===
$ws = $wb.Worksheets.Item("Op2")
$ws.Cells.Item(2, 3).Value = -1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 3).Value = -1
$ws.Cells.Item(3, 4).Value = 5
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(4, 3).Value = -1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 3).Value = -1
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 3).Value = -1
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(7, 3).Value = -1
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 3).Value = -1
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 9).Value = 3
$ws.Cells.Item(9, 3).Value = -1
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 9).Value = 3
$ws.Cells.Item(10, 3).Value = -1
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 9).Value = 2
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(12, 3).Value = -1
$ws.Cells.Item(12, 4).Value = 7
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 3).Value = -1
$ws.Cells.Item(13, 4).Value = 6
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(14, 3).Value = -1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 3).Value = -1
$ws.Cells.Item(15, 4).Value = 6
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(18, 3).Value = -1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(19, 3).Value = -1
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 3).Value = -1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 3).Value = -1
$ws.Cells.Item(21, 4).Value = 5
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 9).Value = 2
$ws.Cells.Item(22, 3).Value = -1
$ws.Cells.Item(22, 4).Value = 28
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 3).Value = -1
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(24, 3).Value = -1
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(25, 3).Value = -1
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(26, 3).Value = -1
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = -1
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 9).Value = 3
$ws.Cells.Item(28, 3).Value = -1
$ws.Cells.Item(28, 4).Value = 18
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = 9
$ws.Cells.Item(29, 3).Value = -1
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = 2
$ws.Cells.Item(30, 3).Value = -1
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 3).Value = -1
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 3).Value = -1
$ws.Cells.Item(32, 4).Value = 6
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 9).Value = 4
$ws.Cells.Item(33, 3).Value = -1
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(34, 3).Value = -1
$ws.Cells.Item(34, 4).Value = 9
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 9).Value = 3

# === Sheet: Barinel ===
$ws = $wb.Worksheets.Item("Barinel")
$ws.Cells.Item(2, 3).Value = -1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 3).Value = -1
$ws.Cells.Item(3, 4).Value = 25
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = 16
$ws.Cells.Item(4, 3).Value = -1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 3).Value = -1
$ws.Cells.Item(5, 4).Value = 8
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 9).Value = 6
$ws.Cells.Item(6, 3).Value = -1
$ws.Cells.Item(6, 4).Value = 11
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 9).Value = 10
$ws.Cells.Item(7, 3).Value = -1
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 9).Value = 28
$ws.Cells.Item(8, 3).Value = -1
$ws.Cells.Item(8, 4).Value = 27
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 9).Value = 13
$ws.Cells.Item(9, 3).Value = -1
$ws.Cells.Item(9, 4).Value = 27
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 9).Value = 13
$ws.Cells.Item(10, 3).Value = -1
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 9).Value = 11
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(12, 3).Value = -1
$ws.Cells.Item(12, 4).Value = 9
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 3).Value = -1
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 9).Value = 6
$ws.Cells.Item(14, 3).Value = -1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 3).Value = -1
$ws.Cells.Item(15, 4).Value = 8
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 9).Value = 6
$ws.Cells.Item(18, 3).Value = -1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(19, 3).Value = -1
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 3).Value = -1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 3).Value = -1
$ws.Cells.Item(21, 4).Value = 11
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 9).Value = 11
$ws.Cells.Item(22, 3).Value = -1
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 3).Value = -1
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(24, 3).Value = -1
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(25, 3).Value = -1
$ws.Cells.Item(25, 4).Value = 26
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 9).Value = 13
$ws.Cells.Item(26, 3).Value = -1
$ws.Cells.Item(26, 4).Value = 6
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = -1
$ws.Cells.Item(27, 4).Value = 26
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 9).Value = 13
$ws.Cells.Item(28, 3).Value = -1
$ws.Cells.Item(28, 4).Value = 22
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = 10
$ws.Cells.Item(29, 3).Value = -1
$ws.Cells.Item(29, 4).Value = 8
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = 12
$ws.Cells.Item(30, 3).Value = -1
$ws.Cells.Item(30, 4).Value = 5
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 3).Value = -1
$ws.Cells.Item(31, 4).Value = 8
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 3).Value = -1
$ws.Cells.Item(32, 4).Value = 8
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 9).Value = 7
$ws.Cells.Item(33, 3).Value = -1
$ws.Cells.Item(33, 4).Value = 9
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(34, 3).Value = -1
$ws.Cells.Item(34, 4).Value = 40
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 9).Value = 13

# === Sheet: Dstar ===
$ws = $wb.Worksheets.Item("Dstar")
$ws.Cells.Item(2, 3).Value = -1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 3).Value = -1
$ws.Cells.Item(3, 4).Value = 5
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(4, 3).Value = -1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 3).Value = -1
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 3).Value = -1
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(7, 3).Value = -1
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 3).Value = -1
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 9).Value = 3
$ws.Cells.Item(9, 3).Value = -1
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 9).Value = 3
$ws.Cells.Item(10, 3).Value = -1
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 9).Value = 2
$ws.Cells.Item(11, 3).Value = -1
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(12, 3).Value = -1
$ws.Cells.Item(12, 4).Value = 7
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(13, 3).Value = -1
$ws.Cells.Item(13, 4).Value = 6
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(14, 3).Value = -1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 3).Value = -1
$ws.Cells.Item(15, 4).Value = 6
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(18, 3).Value = -1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(19, 3).Value = -1
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 3).Value = -1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 3).Value = -1
$ws.Cells.Item(21, 4).Value = 5
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 9).Value = 2
$ws.Cells.Item(22, 3).Value = -1
$ws.Cells.Item(22, 4).Value = 23
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 3).Value = -1
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(24, 3).Value = -1
$ws.Cells.Item(24, 4).Value = 4
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(25, 3).Value = -1
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(26, 3).Value = -1
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 3).Value = -1
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 9).Value = 3
$ws.Cells.Item(28, 3).Value = -1
$ws.Cells.Item(28, 4).Value = 18
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 9).Value = 9
$ws.Cells.Item(29, 3).Value = -1
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 9).Value = 2
$ws.Cells.Item(30, 3).Value = -1
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 3).Value = -1
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 3).Value = -1
$ws.Cells.Item(32, 4).Value = 6
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 9).Value = 4
$ws.Cells.Item(33, 3).Value = -1
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 9).Value = 5
$ws.Cells.Item(34, 3).Value = -1
$ws.Cells.Item(34, 4).Value = 132
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 132
$ws.Cells.Item(34, 7).Value = 132
$ws.Cells.Item(34, 9).Value = 3
